$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "❌ EXPIRED 3768 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 11 days"
$ws.Range("E8").Value = "⚠️ Expires in 11 days"
$ws.Range("E9").Value = "⚠️ Expires in 30 days"
